$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (kept as Text, matching the original inline-string cells)
$updates = @{
    "D2" = "303.69"
    "E2" = "0.49%"
    "G2" = "2"
    "D3" = "31.77"
    "E3" = "-1.59%"
    "G3" = "2"
    "D4" = "5.204"
    "E4" = "0.81%"
    "G4" = "2"
    "D5" = "0.07443"
    "E5" = "-0.94%"
    "G5" = "2"
    "D6" = "2.221"
    "E6" = "31.08%"
    "G6" = "2"
    "D7" = "7.949"
    "E7" = "1.09%"
    "G7" = "2"
    "D8" = "3.858"
    "E8" = "0.32%"
    "G8" = "2"
    "D9" = "0.9183"
    "E9" = "-0.48%"
    "G9" = "2"
    "D10" = "0.1730"
    "E10" = "1.17%"
    "G10" = "2"
    "D11" = "0.07666"
    "E11" = "0.51%"
    "G11" = "2"
    "D12" = "0.08173"
    "E12" = "1.83%"
    "G12" = "2"
    "D13" = "0.02980"
    "E13" = "-0.85%"
    "G13" = "2"
    "D14" = "0.09928"
    "E14" = "0.16%"
    "G14" = "2"
    "D15" = "0.001492"
    "E15" = "0.00%"
    "G15" = "2"
    "D16" = "0.006082"
    "E16" = "-1.92%"
    "G16" = "2"
    "D17" = "3.492"
    "E17" = "0.85%"
    "G17" = "2"
    "D18" = "2.229"
    "E18" = "-0.12%"
    "G18" = "2"
    "D19" = "0.3263"
    "E19" = "-0.95%"
    "G19" = "2"
    "D20" = "0.1343"
    "E20" = "0.52%"
    "G20" = "2"
    "E21" = "1.81%"
    "G21" = "2"
    "D22" = "0.04633"
    "E22" = "0.29%"
    "G22" = "2"
    "D23" = "0.1558"
    "E23" = "0.55%"
    "G23" = "2"
    "D24" = "0.001217"
    "E24" = "-0.04%"
    "G24" = "2"
    "D25" = "0.004504"
    "E25" = "1.51%"
    "G25" = "2"
    "D26" = "0.0001295"
    "E26" = "-7.46%"
    "G26" = "2"
    "D27" = "0.0002729"
    "E27" = "51.57%"
    "G27" = "2"
    "G28" = "2"
    "G29" = "2"
    "G30" = "2"
    "G31" = "2"
    "G32" = "2"
    "G33" = "2"
    "G34" = "2"
    "G35" = "2"
    "G36" = "2"
    "G37" = "2"
    "G38" = "2"
    "D39" = "0.01788"
    "E39" = "7.31%"
    "G39" = "2"
    "D40" = "0.04541"
    "E40" = "-0.50%"
    "G40" = "2"
    "D41" = "0.007318"
    "E41" = "4.78%"
    "G41" = "2"
    "D42" = "0.1357"
    "E42" = "0.99%"
    "G42" = "2"
    "D43" = "0.002162"
    "E43" = "4.98%"
    "G43" = "2"
    "D44" = "0.01092"
    "E44" = "-11.90%"
    "G44" = "2"
    "D45" = "0.00006257"
    "E45" = "3.28%"
    "G45" = "2"
    "E46" = "-57.47%"
    "G46" = "2"
    "D47" = "0.009852"
    "E47" = "-19.50%"
    "G47" = "2"
    "G48" = "2"
    "G49" = "2"
    "G50" = "2"
    "G51" = "2"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"   # force Text storage so values like "303.69" / "0.49%" are not reinterpreted as numbers
    $cell.Value = $updates[$addr]
}
